# Update the AVIV_CheckoutPage sheet's recorded checkout test-data
# (address, phone number, order confirmation number) with the values
# captured from the most recent automated test run.
#
# The values look like plain numbers ("3109847671", "5064293205"), but the
# workbook stores them as text (shared strings), so they are entered with a
# leading apostrophe to force text entry, then the quote-prefix formatting
# that Excel applies to such "number stored as text" cells is cleared by
# resetting the cell style back to Normal - matching how the other text
# cells on the sheet are styled.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AVIV_CheckoutPage")

$ws.Range("D6").Value = "'3109847671"
$ws.Range("D6").Style = "Normal"

$ws.Range("D8").Value = "'5064293205"
$ws.Range("D8").Style = "Normal"

$ws.Range("D13").Value = "ORDER NUMBER: 1038"
$ws.Range("D13").Style = "Normal"
